$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 50) with "NULL" in every cell, A through F.
$ws.Range("A50:F50").Value = "NULL"

# Scroll the view down and move the selection, matching the saved view state.
$ws.Range("C43").Select()
$excel.ActiveWindow.ScrollRow = 31

$ws.Range("C43").Activate()
